$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Update Price (D) and Volume(1h) (E) columns for existing rows
Set-TextValue 'D2' '26.045.76'
Set-TextValue 'E2' '  -0.23%  '
Set-TextValue 'D3' '1.637.50'
Set-TextValue 'E3' '  -1.78%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '213.49'
Set-TextValue 'E5' '  +1.96%  '
Set-TextValue 'D6' '0.5235'
Set-TextValue 'E6' '  -0.28%  '
Set-TextValue 'E7' '  -0.10%  '
Set-TextValue 'E8' '  -0.70%  '
Set-TextValue 'D9' '0.06288'
Set-TextValue 'E10' '  -2.43%  '
Set-TextValue 'D11' '0.07675'
Set-TextValue 'E11' '  +1.96%  '
Set-TextValue 'D12' '1.641.37'
Set-TextValue 'E12' '  -1.52%  '
Set-TextValue 'D13' '4.392'
Set-TextValue 'E13' '  -0.74%  '
Set-TextValue 'D14' '1.858.07'
Set-TextValue 'E14' '  -1.96%  '
Set-TextValue 'D15' '0.5503'
Set-TextValue 'E15' '  +0.15%  '
Set-TextValue 'D16' '0.0₅8175'
Set-TextValue 'E16' '  +2.92%  '
Set-TextValue 'D17' '64.77'
Set-TextValue 'E17' '  -2.35%  '
Set-TextValue 'D18' '26.042.70'
Set-TextValue 'E18' '  -0.32%  '
Set-TextValue 'E19' '  -0.06%  '
Set-TextValue 'D20' '4.673'
Set-TextValue 'E20' '  -0.62%  '
Set-TextValue 'D21' '187.93'
Set-TextValue 'E21' '  +0.84%  '
Set-TextValue 'D22' '10.16'
Set-TextValue 'E22' '  -0.78%  '
Set-TextValue 'D23' '6.141'
Set-TextValue 'E23' '  -0.38%  '
Set-TextValue 'E24' '  -0.07%  '
Set-TextValue 'D25' '145.19'
Set-TextValue 'E25' '  -2.96%  '
Set-TextValue 'E26' '  -2.79%  '
Set-TextValue 'D27' '7.385'
Set-TextValue 'E27' '  -0.86%  '
Set-TextValue 'E28' '  -0.67%  '
Set-TextValue 'D29' '1.377'
Set-TextValue 'E29' '  +1.99%  '
Set-TextValue 'D30' '0.05947'
Set-TextValue 'E30' '  -6.64%  '
Set-TextValue 'D31' '1.253'
Set-TextValue 'E31' '  -1.56%  '
Set-TextValue 'E32' '  -1.94%  '
Set-TextValue 'D33' '3.392'
Set-TextValue 'E33' '  -0.47%  '
Set-TextValue 'D34' '1.640'
Set-TextValue 'E34' '  +0.30%  '
Set-TextValue 'D35' '0.9797'
Set-TextValue 'E35' '  -2.22%  '
Set-TextValue 'E36' '  -0.41%  '
Set-TextValue 'D37' '2.761'
Set-TextValue 'E37' '  +1.15%  '
Set-TextValue 'D38' '0.5626'
Set-TextValue 'E38' '  -6.27%  '
Set-TextValue 'E39' '  -0.15%  '
Set-TextValue 'D40' '0.8486'
Set-TextValue 'E40' '  -2.49%  '
Set-TextValue 'D42' '5.679'
Set-TextValue 'E42' '  -6.84%  '
Set-TextValue 'D43' '1.025.71'
Set-TextValue 'E43' '  -7.37%  '
Set-TextValue 'D44' '100.06'
Set-TextValue 'E44' '  +0.28%  '
Set-TextValue 'D45' '1.785.72'
Set-TextValue 'E45' '  -1.76%  '

# Rows 46-51: a new row (BabyDogeCoin) was inserted at 46, shifting Aave..Mantle
# down by one and dropping Aptos off the bottom of the list. Apply as direct
# content overwrites (A column rank index cells are untouched).
Set-TextValue 'B46' 'BabyDogeCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.0₈108'
Set-TextValue 'E46' '  +0.52%  '
Set-TextValue 'B47' 'Aave'
Set-TextValue 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '55.66'
Set-TextValue 'E47' '  +0.77%  '
Set-TextValue 'B48' 'Frax'
Set-TextValue 'C48' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D48' '0.9997'
Set-TextValue 'E48' '  -0.31%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '8.045'
Set-TextValue 'E49' '  +0.15%  '
Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.05141'
Set-TextValue 'E50' '  -1.64%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.4217'
Set-TextValue 'E51' '  -0.65%  '
